$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Default Values")

# Remove the "modules__school__schoolName" column (column G) from the
# "Default Values" sheet.
$ws.Columns.Item(7).Delete()

# Re-assert the header formatting on the remaining header cells so the
# engine collapses them back onto the shared bold style instead of
# leaving them pinned to a leftover duplicate style slot.
$ws.Range("F1:H1").Font.Bold = $true

# Move selection to match the target workbook state.
$ws.Activate()
$ws.Range("F6").Select()

# The "Connectors" sheet's header row was pinned to a redundant duplicate
# style slot as well; re-assert it so it collapses back onto the shared
# (non-bold) default style.
$ws1 = $wb.Worksheets.Item("Connectors")
$ws1.Range("A1:H1").Font.Bold = $false
